$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 7458
$ws.Range("J3").Value = 7842
$ws.Range("E4").Value = 1971
$ws.Range("G4").Value = 1437
$ws.Range("J4").Value = 1707
$ws.Range("J5").Value = 613
$ws.Range("J6").Value = 10712
$ws.Range("E7").Value = 25266
$ws.Range("G7").Value = 24112
$ws.Range("J7").Value = 28332

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 77
$ws.Range("J3").Value = 56
$ws.Range("J7").Value = 426

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 472
$ws.Range("J3").Value = 515
$ws.Range("J6").Value = 665
$ws.Range("J7").Value = 1790

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 169
$ws.Range("J6").Value = 154
$ws.Range("J7").Value = 571

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 428
$ws.Range("J4").Value = 56
$ws.Range("J6").Value = 456
$ws.Range("J7").Value = 1285

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 143
$ws.Range("J7").Value = 405

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 263
$ws.Range("J3").Value = 290
$ws.Range("J7").Value = 869

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 207
$ws.Range("J3").Value = 190
$ws.Range("J4").Value = 31
$ws.Range("J6").Value = 266
$ws.Range("J7").Value = 713

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 224
$ws.Range("J6").Value = 217
$ws.Range("J7").Value = 809
$ws.Range("J8").Value = 1790
$ws.Range("J11").Value = 507
$ws.Range("J18").Value = 229
$ws.Range("J23").Value = 261
$ws.Range("J24").Value = 99
$ws.Range("J25").Value = 147
$ws.Range("J26").Value = 54
$ws.Range("J27").Value = 170
$ws.Range("J29").Value = 1506
$ws.Range("J32").Value = 48
$ws.Range("J33").Value = 1285
$ws.Range("J34").Value = 131
$ws.Range("J36").Value = 384
$ws.Range("J37").Value = 869
$ws.Range("J41").Value = 215
$ws.Range("J42").Value = 1206
$ws.Range("J44").Value = 223
$ws.Range("J50").Value = 172
$ws.Range("J52").Value = 721
$ws.Range("J53").Value = 426
$ws.Range("J54").Value = 562
$ws.Range("J55").Value = 444
$ws.Range("G63").Value = 256
$ws.Range("J63").Value = 83
$ws.Range("J65").Value = 713
$ws.Range("J67").Value = 1033
$ws.Range("J72").Value = 108
$ws.Range("J76").Value = 401
$ws.Range("J77").Value = 200
$ws.Range("J78").Value = 328
$ws.Range("J79").Value = 774
$ws.Range("J83").Value = 571
$ws.Range("E85").Value = 974
$ws.Range("J85").Value = 1167
$ws.Range("J88").Value = 301
$ws.Range("J89").Value = 355
$ws.Range("J90").Value = 297
$ws.Range("J95").Value = 405
$ws.Range("E101").Value = 25266
$ws.Range("G101").Value = 24112
$ws.Range("J101").Value = 28332

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J6").Value = 287
$ws.Range("J7").Value = 1033

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 261
$ws.Range("J7").Value = 562

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J3").Value = 530
$ws.Range("J7").Value = 1506

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 67
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 90
$ws.Range("J7").Value = 401

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 65
$ws.Range("J3").Value = 51
$ws.Range("J7").Value = 217

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J4").Value = 11
$ws.Range("J6").Value = 130
$ws.Range("J7").Value = 215

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 252
$ws.Range("J6").Value = 640
$ws.Range("J7").Value = 1206

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 328

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 444

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J2").Value = 73
$ws.Range("J7").Value = 261

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 221
$ws.Range("J3").Value = 258
$ws.Range("J7").Value = 774

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J6").Value = 107
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 384

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 255
$ws.Range("J3").Value = 242
$ws.Range("J5").Value = 22
$ws.Range("J6").Value = 257
$ws.Range("J7").Value = 809

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 147

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J2").Value = 44
$ws.Range("J7").Value = 172

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J4").Value = 30
$ws.Range("J6").Value = 242
$ws.Range("J7").Value = 507

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 224

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 160
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 101
$ws.Range("J7").Value = 355

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 62
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 75
$ws.Range("J4").Value = 32

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 312
$ws.Range("J3").Value = 422
$ws.Range("E4").Value = 59
$ws.Range("J4").Value = 72
$ws.Range("E7").Value = 974
$ws.Range("J7").Value = 1167

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 75
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 200

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 170
$ws.Range("J6").Value = 311
$ws.Range("J7").Value = 721
